# Generate Report for Archive
# - Update the localization status text from "Ready for handoff" to
#   "In Translation" on every sheet that reports it (Overview summary
#   columns for each language, plus each language detail sheet's
#   "Status" column).
# - Narrow the affected "Status" columns to match the new (shorter)
#   text's auto-fit width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: column E = zh-cn status, column F = de-de status
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Per-language detail sheets: column C = Status
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Re-fit the now-narrower Status columns
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
